$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.457.48"
$ws.Range("E2").Value = "  +0.72%  "
$ws.Range("D3").Value = "1.617.12"
$ws.Range("E3").Value = "  +1.45%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "'212.51"
$ws.Range("E5").Value = "  -0.26%  "
$ws.Range("D6").Value = "'0.499"
$ws.Range("E6").Value = "  -0.21%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("E8").Value = "  -0.11%  "
$ws.Range("E9").Value = "  +0.11%  "
$ws.Range("D10").Value = "'19.25"
$ws.Range("E10").Value = "  +1.36%  "
$ws.Range("E11").Value = "  -0.26%  "
$ws.Range("D12").Value = "1.844.84"
$ws.Range("E12").Value = "  +1.43%  "
$ws.Range("D13").Value = "1.610.60"
$ws.Range("E13").Value = "  +1.36%  "
$ws.Range("E14").Value = "  +0.42%  "
$ws.Range("E15").Value = "  +0.10%  "
$ws.Range("D16").Value = "'63.75"
$ws.Range("E16").Value = "  -0.08%  "
$ws.Range("D17").Value = "'235.05"
$ws.Range("E17").Value = "  +9.02%  "
$ws.Range("D18").Value = "26.460.90"
$ws.Range("E18").Value = "  +0.76%  "
$ws.Range("D19").Value = "'7.72"
$ws.Range("E19").Value = "  +5.23%  "
$ws.Range("D20").Value = "0.0₃0725"
$ws.Range("E20").Value = "  +0.38%  "
$ws.Range("D21").Value = "'0.999"
$ws.Range("E21").Value = "  -0.14%  "
$ws.Range("E22").Value = "  -0.30%  "
$ws.Range("E23").Value = "  +4.48%  "
$ws.Range("E24").Value = "  +0.32%  "
$ws.Range("D25").Value = "'146.97"
$ws.Range("E25").Value = "  +1.46%  "
$ws.Range("E26").Value = "  -0.13%  "
$ws.Range("D27").Value = "'7.00"
$ws.Range("E27").Value = "  +0.62%  "
$ws.Range("E28").Value = "  +0.14%  "
$ws.Range("D29").Value = "'15.53"
$ws.Range("E29").Value = "  +2.79%  "
$ws.Range("E30").Value = "  +0.99%  "
$ws.Range("E31").Value = "  -0.27%  "
$ws.Range("D32").Value = "1.510.35"
$ws.Range("E32").Value = "  +6.58%  "
$ws.Range("E33").Value = "  +1.69%  "
$ws.Range("D34").Value = "'2.96"
$ws.Range("E34").Value = "  -0.22%  "
$ws.Range("D35").Value = "'1.53"
$ws.Range("E35").Value = "  +4.57%  "
$ws.Range("E36").Value = "  -0.21%  "
$ws.Range("D37").Value = "'0.564"
$ws.Range("E37").Value = "  -1.95%  "
$ws.Range("E38").Value = "  +0.00%  "
$ws.Range("D39").Value = "'0.829"
$ws.Range("E39").Value = "  +0.72%  "
$ws.Range("D40").Value = "'5.90"
$ws.Range("E40").Value = "  +1.83%  "
$ws.Range("E41").Value = "  -0.06%  "
$ws.Range("E42").Value = "  +1.27%  "
$ws.Range("D43").Value = "1.757.25"
$ws.Range("E43").Value = "  +1.50%  "
$ws.Range("D44").Value = "'0.760"
$ws.Range("E44").Value = "  -0.16%  "
$ws.Range("D45").Value = "'61.52"
$ws.Range("E45").Value = "  +0.99%  "
$ws.Range("D46").Value = "'0.908"
$ws.Range("E46").Value = "  -2.03%  "
$ws.Range("D47").Value = "'89.86"
$ws.Range("E47").Value = "  +3.13%  "
$ws.Range("D48").Value = "'1.49"
$ws.Range("E48").Value = "  +1.19%  "
$ws.Range("E49").Value = "  +0.10%  "
$ws.Range("D50").Value = "'0.0962"
$ws.Range("E50").Value = "  +0.95%  "
$ws.Range("E51").Value = "  +1.16%  "
